# aggiornamento fino a 6/03
# Appends three new daily rows (245-247) to the Castelvetro report sheet,
# continuing the date series in column A and copying the date-column style
# from the last existing row so formatting stays consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 244

$newRecords = @(
    @(44319, 0, 18, 159.4472495349455),
    @(44320, 1, 12, 106.2981663566303),
    @(44321, 2, 14, 124.0145274160687)
)

$rowIdx = $lastRow + 1
foreach ($rec in $newRecords) {
    # Copy the style of the date cell from the previous row so the new
    # date cell keeps the same formatting (border/alignment/date format).
    $ws.Cells.Item($lastRow, 1).Copy($ws.Cells.Item($rowIdx, 1))

    $ws.Cells.Item($rowIdx, 1).Value = $rec[0]
    $ws.Cells.Item($rowIdx, 2).Value = $rec[1]
    $ws.Cells.Item($rowIdx, 3).Value = $rec[2]
    $ws.Cells.Item($rowIdx, 4).Value = $rec[3]

    $rowIdx++
}
